$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (values Excel will not misinterpret as numbers)
$ws.Range("D2").Value = '61.299.52'
$ws.Range("E2").Value = '  -1.11%  '
$ws.Range("D3").Value = '2.432.77'
$ws.Range("E3").Value = '  -0.20%  '
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("E6").Value = '  -1.40%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  +0.49%  '
$ws.Range("D9").Value = '2.418.33'
$ws.Range("E10").Value = '  +2.15%  '
$ws.Range("E11").Value = '  +1.71%  '
$ws.Range("E12").Value = '  -0.75%  '
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("E14").Value = '  -0.51%  '
$ws.Range("D15").Value = '2.893.03'
$ws.Range("E15").Value = '  +0.61%  '
$ws.Range("E16").Value = '  -0.84%  '
$ws.Range("D17").Value = '61.252.42'
$ws.Range("E17").Value = '  -1.21%  '
$ws.Range("D18").Value = '2.431.85'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("E19").Value = '  -2.72%  '
$ws.Range("E20").Value = '  +2.58%  '
$ws.Range("E21").Value = '  -1.29%  '
$ws.Range("E22").Value = '  -0.96%  '
$ws.Range("E23").Value = '  +1.51%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("E25").Value = '  -1.73%  '
$ws.Range("E26").Value = '  -0.77%  '
$ws.Range("E27").Value = '  -4.63%  '
$ws.Range("E28").Value = '  -7.15%  '
$ws.Range("E29").Value = '  +0.48%  '
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("D31").Value = '0.0₃0915'
$ws.Range("E31").Value = '  -2.74%  '
$ws.Range("E32").Value = '  -1.05%  '
$ws.Range("E33").Value = '  -5.04%  '
$ws.Range("E34").Value = '  -0.87%  '
$ws.Range("E35").Value = '  -6.94%  '
$ws.Range("E36").Value = '  +0.27%  '
$ws.Range("E37").Value = '  -5.57%  '
$ws.Range("E38").Value = '  +1.20%  '
$ws.Range("E39").Value = '  -1.20%  '
$ws.Range("E40").Value = '  -2.96%  '
$ws.Range("E41").Value = '  +0.38%  '
$ws.Range("E42").Value = '  -2.11%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("E44").Value = '  -2.55%  '
$ws.Range("E45").Value = '  -4.98%  '
$ws.Range("E46").Value = '  -3.93%  '
$ws.Range("D47").Value = '0.0₆0289'
$ws.Range("E47").Value = '  +24.99%  '
$ws.Range("E48").Value = '  -0.73%  '
$ws.Range("E49").Value = '  -2.06%  '
$ws.Range("E50").Value = '  -0.65%  '
$ws.Range("E51").Value = '  -2.88%  '

# Numeric-looking price updates that must stay as literal text (e.g. "8.90" keeps trailing zero)
# Mark the cell as Text format first so Excel keeps the exact string instead of coercing to a float
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.997'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '574.66'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.61'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.13'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.339'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.16'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.59'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.24'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '324.01'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '65.03'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.90'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '571.63'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.89'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.34'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.133'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.61'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '151.87'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.369'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.31'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.11'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.69'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.66'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.35'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '141.00'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.592'
